# The commit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml  ("Office Theme" colors) <-> ppt/theme/theme2.xml ("Integral" colors)
# theme2.xml is the theme actually wired to the slide master / presentation (the one
# PowerPoint's object model exposes and renders), so the visible, user-facing effect of
# the swap is that the deck's applied theme colors change from the "Integral" palette to
# the classic "Office Theme" palette. We reproduce that effect by rewriting the active
# theme's 12-color scheme (ThemeColorScheme, in dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink
# order) to the Office Theme values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Office Theme color scheme, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5,
# 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's RGB color value is packed as 0x00BBGGRR (BGR byte order).
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $cs.Item($i).RGB = $bgr
}
